$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45204 -> 45205) for every data row (rows 2 through 458).
$ws.Range("C2:C458").Value = 45205
